$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 848.1429000000001  # H38: was 575.25
$ws.Cells.Item(38, 9).Value = 67.71429000000001  # I38: was 132.66667
$ws.Cells.Item(38, 10).Value = 1628.5714  # J38: was 1903
$ws.Cells.Item(38, 11).Value = 203.14287  # K38: was 398.00001
$ws.Cells.Item(38, 12).Value = 4885.7142  # L38: was 5709
$ws.Cells.Item(38, 13).Value = 168.85713  # M38: was -26.00001000000003
$ws.Cells.Item(38, 14).Value = -5629.7142  # N38: was -6453

$ws.Cells.Item(64, 8).Value = 4100  # H64: was 3990
$ws.Cells.Item(64, 10).Value = 4342.857  # J64: was 4175
$ws.Cells.Item(64, 12).Value = 4342.857  # L64: was 4175
$ws.Cells.Item(64, 14).Value = -4838.857  # N64: was -4671

$ws.Cells.Item(67, 8).Value = 4100  # H67: was 3990
$ws.Cells.Item(67, 10).Value = 4342.857  # J67: was 4175
$ws.Cells.Item(67, 12).Value = 4342.857  # L67: was 4175
$ws.Cells.Item(67, 14).Value = -6058.857  # N67: was -5891

$ws.Cells.Item(116, 8).Value = 5566.3335  # H116: was 5655.091
$ws.Cells.Item(116, 9).Value = 3333.3333  # I116: was 3133.3333
$ws.Cells.Item(116, 10).Value = 6310.6665  # J116: was 6600.75
$ws.Cells.Item(116, 11).Value = 3333.3333  # K116: was 3133.3333
$ws.Cells.Item(116, 12).Value = 6310.6665  # L116: was 6600.75
$ws.Cells.Item(116, 13).Value = 108.6667000000002  # M116: was 308.6667000000002
$ws.Cells.Item(116, 14).Value = -13194.6665  # N116: was -13484.75

$ws.Cells.Item(132, 8).Value = 3077.8708  # H132: was 2721.639
$ws.Cells.Item(132, 9).Value = 3300.2856  # I132: was 2808.0293
$ws.Cells.Item(132, 10).Value = 1002  # J132: was 1253
$ws.Cells.Item(132, 11).Value = 9900.856800000001  # K132: was 8424.0879
$ws.Cells.Item(132, 12).Value = 3006  # L132: was 3759
$ws.Cells.Item(132, 13).Value = -7370.856800000001  # M132: was -5894.0879
$ws.Cells.Item(132, 14).Value = -8066  # N132: was -8819

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(24, 8).Value = 15499  # H24: was 15249.5
$ws.Cells.Item(24, 10).Value = 15499  # J24: was 15249.5
$ws.Cells.Item(24, 12).Value = 15499  # L24: was 15249.5
$ws.Cells.Item(24, 14).Value = -16247  # N24: was -15997.5

$ws.Cells.Item(32, 8).Value = 5025.0493  # H32: was 4632.5317
$ws.Cells.Item(32, 9).Value = 4335.5557  # I32: was 3814.3252
$ws.Cells.Item(32, 10).Value = 10541  # J32: was 10806.272
$ws.Cells.Item(32, 11).Value = 4335.5557  # K32: was 3814.3252
$ws.Cells.Item(32, 12).Value = 10541  # L32: was 10806.272
$ws.Cells.Item(32, 13).Value = -4048.5557  # M32: was -3527.3252
$ws.Cells.Item(32, 14).Value = -11115  # N32: was -11380.272

$ws.Cells.Item(45, 8).Value = 2395.6428  # H45: was 2411.7144
$ws.Cells.Item(45, 9).Value = 1900.6154  # I45: was 1829.1428
$ws.Cells.Item(45, 10).Value = 2824.6667  # J45: was 2994.2856
$ws.Cells.Item(45, 11).Value = 1900.6154  # K45: was 1829.1428
$ws.Cells.Item(45, 12).Value = 2824.6667  # L45: was 2994.2856
$ws.Cells.Item(45, 13).Value = -1523.6154  # M45: was -1452.1428
$ws.Cells.Item(45, 14).Value = -3578.6667  # N45: was -3748.2856

$ws.Cells.Item(74, 8).Value = 22728588  # H74: was 23257142
$ws.Cells.Item(74, 10).Value = 3310.3333  # J74: was 3547.818
$ws.Cells.Item(74, 12).Value = 3310.3333  # L74: was 3547.818
$ws.Cells.Item(74, 14).Value = -5058.3333  # N74: was -5295.818

$ws.Cells.Item(77, 8).Value = 22728588  # H77: was 23257142
$ws.Cells.Item(77, 10).Value = 3310.3333  # J77: was 3547.818
$ws.Cells.Item(77, 12).Value = 16551.6665  # L77: was 17739.09
$ws.Cells.Item(77, 14).Value = -25287.6665  # N77: was -26475.09

$ws.Cells.Item(100, 8).Value = 15499  # H100: was 15249.5
$ws.Cells.Item(100, 10).Value = 15499  # J100: was 15249.5
$ws.Cells.Item(100, 12).Value = 15499  # L100: was 15249.5
$ws.Cells.Item(100, 14).Value = -17663  # N100: was -17413.5

$ws.Cells.Item(122, 8).Value = 1570  # H122: was 1481.7391
$ws.Cells.Item(122, 9).Value = 1197.6111  # I122: was 1133.35
$ws.Cells.Item(122, 11).Value = 3592.8333  # K122: was 3400.05
$ws.Cells.Item(122, 13).Value = -1142.8333  # M122: was -950.0499999999997

$ws.Cells.Item(132, 8).Value = 10075.55  # H132: was 9048.940000000001
$ws.Cells.Item(132, 9).Value = 1551.2391  # I132: was 1379.3019
$ws.Cells.Item(132, 11).Value = 4653.7173  # K132: was 4137.905699999999
$ws.Cells.Item(132, 13).Value = -2123.7173  # M132: was -1607.905699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3098.2092  # H134: was 3155.024
$ws.Cells.Item(134, 9).Value = 3111.8057  # I134: was 3180.3713
$ws.Cells.Item(134, 11).Value = 9335.417099999999  # K134: was 9541.1139
$ws.Cells.Item(134, 13).Value = -6800.417099999999  # M134: was -7006.1139

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3327.2046  # H31: was 3551.5
$ws.Cells.Item(31, 9).Value = 2663.1304  # I31: was 2724.9092
$ws.Cells.Item(31, 10).Value = 4054.524  # J31: was 4688.0625
$ws.Cells.Item(31, 11).Value = 2663.1304  # K31: was 2724.9092
$ws.Cells.Item(31, 12).Value = 4054.524  # L31: was 4688.0625
$ws.Cells.Item(31, 13).Value = -2368.1304  # M31: was -2429.9092
$ws.Cells.Item(31, 14).Value = -4644.523999999999  # N31: was -5278.0625

$ws.Cells.Item(34, 8).Value = 3327.2046  # H34: was 3551.5
$ws.Cells.Item(34, 9).Value = 2663.1304  # I34: was 2724.9092
$ws.Cells.Item(34, 10).Value = 4054.524  # J34: was 4688.0625
$ws.Cells.Item(34, 11).Value = 2663.1304  # K34: was 2724.9092
$ws.Cells.Item(34, 12).Value = 4054.524  # L34: was 4688.0625
$ws.Cells.Item(34, 13).Value = -2461.1304  # M34: was -2522.9092
$ws.Cells.Item(34, 14).Value = -4458.523999999999  # N34: was -5092.0625

$ws.Cells.Item(36, 8).Value = 2248  # H36: was 100048
$ws.Cells.Item(36, 9).Value = 2248  # I36: was 100048
$ws.Cells.Item(36, 11).Value = 2248  # K36: was 100048
$ws.Cells.Item(36, 13).Value = -1860  # M36: was -99660

$ws.Cells.Item(40, 8).Value = 2248  # H40: was 100048
$ws.Cells.Item(40, 9).Value = 2248  # I40: was 100048
$ws.Cells.Item(40, 11).Value = 2248  # K40: was 100048
$ws.Cells.Item(40, 13).Value = -2088  # M40: was -99888

$ws.Cells.Item(132, 8).Value = 2385.7188  # H132: was 2455.8386
$ws.Cells.Item(132, 9).Value = 1572.12  # I132: was 1628.7916
$ws.Cells.Item(132, 11).Value = 4716.36  # K132: was 4886.3748
$ws.Cells.Item(132, 13).Value = -2186.36  # M132: was -2356.3748

$ws.Cells.Item(134, 8).Value = 1147.7097  # H134: was 1289
$ws.Cells.Item(134, 9).Value = 993  # I134: was 1107.6
$ws.Cells.Item(134, 10).Value = 1525.8889  # J134: was 1629.125
$ws.Cells.Item(134, 11).Value = 2979  # K134: was 3322.8
$ws.Cells.Item(134, 12).Value = 4577.6667  # L134: was 4887.375
$ws.Cells.Item(134, 13).Value = -444  # M134: was -787.7999999999997
$ws.Cells.Item(134, 14).Value = -9647.6667  # N134: was -9957.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1268.6666  # H5: was 1535.0555
$ws.Cells.Item(5, 9).Value = 672.125  # I5: was 687.7
$ws.Cells.Item(5, 10).Value = 1950.4286  # J5: was 2594.25
$ws.Cells.Item(5, 11).Value = 2016.375  # K5: was 2063.1
$ws.Cells.Item(5, 12).Value = 5851.2858  # L5: was 7782.75
$ws.Cells.Item(5, 13).Value = -1904.375  # M5: was -1951.1
$ws.Cells.Item(5, 14).Value = -6075.2858  # N5: was -8006.75

$ws.Cells.Item(23, 8).Value = 615  # H23: was 354.72726
$ws.Cells.Item(23, 9).Value = 0  # I23: was 11
$ws.Cells.Item(23, 10).Value = 615  # J23: was 431.1111
$ws.Cells.Item(23, 11).Value = 0  # K23: was 33
$ws.Cells.Item(23, 12).Value = 1845  # L23: was 1293.3333
$ws.Cells.Item(23, 13).ClearContents()  # M23: was 202
$ws.Cells.Item(23, 14).Value = -2315  # N23: was -1763.3333

$ws.Cells.Item(41, 8).Value = 781.3333  # H41: was 519.8
$ws.Cells.Item(41, 9).Value = 422  # I41: was 274.75
$ws.Cells.Item(41, 11).Value = 1266  # K41: was 824.25
$ws.Cells.Item(41, 13).Value = -928  # M41: was -486.25

$ws.Cells.Item(131, 8).Value = 700.45  # H131: was 700.8099999999999
$ws.Cells.Item(131, 10).Value = 700.45  # J131: was 700.8099999999999
$ws.Cells.Item(131, 12).Value = 2101.35  # L131: was 2102.43
$ws.Cells.Item(131, 14).Value = -12181.35  # N131: was -12182.43

$ws.Cells.Item(135, 8).Value = 1268.6666  # H135: was 1535.0555
$ws.Cells.Item(135, 9).Value = 672.125  # I135: was 687.7
$ws.Cells.Item(135, 10).Value = 1950.4286  # J135: was 2594.25
$ws.Cells.Item(135, 11).Value = 6049.125  # K135: was 6189.3
$ws.Cells.Item(135, 12).Value = 17553.8574  # L135: was 23348.25
$ws.Cells.Item(135, 13).Value = -3514.125  # M135: was -3654.3
$ws.Cells.Item(135, 14).Value = -22623.8574  # N135: was -28418.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 74076600  # H122: was 55557652
$ws.Cells.Item(122, 9).Value = 27779148  # I122: was 20834566
$ws.Cells.Item(122, 10).Value = 166671500  # J122: was 125003820
$ws.Cells.Item(122, 11).Value = 83337444  # K122: was 62503698
$ws.Cells.Item(122, 12).Value = 500014500  # L122: was 375011460
$ws.Cells.Item(122, 13).Value = -83334994  # M122: was -62501248
$ws.Cells.Item(122, 14).Value = -500019400  # N122: was -375016360

$ws.Cells.Item(132, 8).Value = 13997.5  # H132: was 12068.038
$ws.Cells.Item(132, 9).Value = 2547.3  # I132: was 2352.8635
$ws.Cells.Item(132, 10).Value = 128499.5  # J132: was 65501.5
$ws.Cells.Item(132, 11).Value = 7641.900000000001  # K132: was 7058.5905
$ws.Cells.Item(132, 12).Value = 385498.5  # L132: was 196504.5
$ws.Cells.Item(132, 13).Value = -5111.900000000001  # M132: was -4528.5905

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 566.8333  # H16: was 604.7273
$ws.Cells.Item(16, 9).Value = 566.8333  # I16: was 604.7273
$ws.Cells.Item(16, 11).Value = 566.8333  # K16: was 604.7273
$ws.Cells.Item(16, 13).Value = -396.8333  # M16: was -434.7273

$ws.Cells.Item(68, 8).Value = 2766  # H68: was 2999
$ws.Cells.Item(68, 10).Value = 2766  # J68: was 2999
$ws.Cells.Item(68, 12).Value = 2766  # L68: was 2999
$ws.Cells.Item(68, 14).Value = -4264  # N68: was -4497

$ws.Cells.Item(71, 8).Value = 2766  # H71: was 2999
$ws.Cells.Item(71, 10).Value = 2766  # J71: was 2999
$ws.Cells.Item(71, 12).Value = 13830  # L71: was 14995
$ws.Cells.Item(71, 14).Value = -21318  # N71: was -22483

$ws.Cells.Item(82, 8).Value = 747.0909  # H82: was 709.8333
$ws.Cells.Item(82, 9).Value = 721.8  # I82: was 683.4545000000001
$ws.Cells.Item(82, 11).Value = 721.8  # K82: was 683.4545000000001
$ws.Cells.Item(82, 13).Value = -360.8  # M82: was -322.4545000000001

$ws.Cells.Item(85, 8).Value = 747.0909  # H85: was 709.8333
$ws.Cells.Item(85, 9).Value = 721.8  # I85: was 683.4545000000001
$ws.Cells.Item(85, 11).Value = 721.8  # K85: was 683.4545000000001
$ws.Cells.Item(85, 13).Value = 526.2  # M85: was 564.5454999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4500.25  # H62: was 4667
$ws.Cells.Item(62, 9).Value = 3666.3333  # I62: was 3750
$ws.Cells.Item(62, 10).Value = 5000.6  # J62: was 4929
$ws.Cells.Item(62, 11).Value = 3666.3333  # K62: was 3750
$ws.Cells.Item(62, 12).Value = 5000.6  # L62: was 4929
$ws.Cells.Item(62, 13).Value = -3042.3333  # M62: was -3126
$ws.Cells.Item(62, 14).Value = -6248.6  # N62: was -6177

$ws.Cells.Item(65, 8).Value = 4500.25  # H65: was 4667
$ws.Cells.Item(65, 9).Value = 3666.3333  # I65: was 3750
$ws.Cells.Item(65, 10).Value = 5000.6  # J65: was 4929
$ws.Cells.Item(65, 11).Value = 18331.6665  # K65: was 18750
$ws.Cells.Item(65, 12).Value = 25003  # L65: was 24645
$ws.Cells.Item(65, 13).Value = -15211.6665  # M65: was -15630
$ws.Cells.Item(65, 14).Value = -31243  # N65: was -30885

$ws.Cells.Item(136, 8).Value = 26471264  # H136: was 27902080
$ws.Cells.Item(136, 9).Value = 39703696  # I136: was 43012240
$ws.Cells.Item(136, 11).Value = 119111088  # K136: was 129036720
$ws.Cells.Item(136, 13).Value = -119108538  # M136: was -129034170
